$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record above the current row 108, pushing the
# existing rows 108-136 down to 109-137 (dimension grows from R136 to R137).
$ws.Rows.Item(108).Insert()

# Populate the new row 108 with the latest weekly observation. It mirrors
# the (now shifted-down) row 109 data except for a new date and volume.
$ws.Cells.Item(108, 1).Value = 8
$ws.Cells.Item(108, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(108, 3).Value = "Coquimbo"
$ws.Cells.Item(108, 4).Value = 44663
$ws.Cells.Item(108, 5).Value = 4
$ws.Cells.Item(108, 6).Value = 100112040
$ws.Cells.Item(108, 7).Value = "Cilantro"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 3000
$ws.Cells.Item(108, 11).Value = 2000
$ws.Cells.Item(108, 12).Value = 2500
$ws.Cells.Item(108, 13).Value = 2250
$ws.Cells.Item(108, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(108, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(108, 16).Value = 1500
$ws.Cells.Item(108, 17).Value = 1.5
$ws.Cells.Item(108, 18).Value = "Hortaliza"
